# Update "想去人数" (column F) figures on the 展览 and 全部类型 sheets
# to reflect newly scraped counts (gh-pages data refresh at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")

$sheet1Updates = @{
    3  = 2813
    4  = 44
    5  = 3096
    6  = 1941
    7  = 252
    9  = 2608
    10 = 588
    11 = 289
    12 = 21
    15 = 136
    16 = 9802
    18 = 20
    19 = 18
    20 = 7742
    21 = 12296
    22 = 126
    24 = 257
    25 = 389
    27 = 2764
    29 = 222
    30 = 2792
    31 = 1281
    33 = 64
    35 = 4583
    36 = 1238
    37 = 46
    39 = 67
    40 = 603
}

foreach ($row in $sheet1Updates.Keys) {
    $ws1.Range("F$row").Value = $sheet1Updates[$row]
}

# --- Sheet "全部类型" (sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")

$sheet4Updates = @{
    4  = 2813
    5  = 44
    7  = 3096
    9  = 1941
    11 = 252
    12 = 2608
    14 = 588
    15 = 290
    16 = 21
    18 = 136
    19 = 9802
    20 = 20
    21 = 18
    22 = 7742
    23 = 12296
    25 = 126
    26 = 257
    30 = 2764
    33 = 222
    35 = 64
    37 = 4583
    45 = 603
}

foreach ($row in $sheet4Updates.Keys) {
    $ws4.Range("F$row").Value = $sheet4Updates[$row]
}
